# Lecture 4 deck: "Added fourth practice + minor edits"
#
# Main content change: the second example on the reduce-2 slide
# (list-concatenation reduce) showed the wrong result "210" (copy/paste
# leftover from the multiplication example on the previous slide) -
# fix it to the actual result of reduce(lambda res, x: res + x,
# [[1, 2], [3, 4, 5]], []), i.e. [1, 2, 3, 4, 5].
#
# Plus a handful of 1-EMU nudges to shape position/size left behind by
# the edit (PowerPoint re-lays-out text boxes on edit/save).
#
# NOTE: Shape.Left/Top/Width/Height are exposed as single-precision
# (32-bit) floats in points, while OOXML stores EMU (1 pt = 12700 EMU)
# as integers. Converting target_emu -> points -> (float32) -> back to
# EMU truncates instead of rounding, which would silently knock values
# off by 1 EMU. Adding half an EMU (in points) before the float32 cast
# compensates for the truncation so the written <a:off>/<a:ext> values
# land exactly on the intended EMU.

function Set-ShapeEmu($Shape, $Left, $Top, $Width, $Height) {
    $halfEmuPt = 0.5 / 12700.0
    $Shape.Left = ($Left / 12700.0) + $halfEmuPt
    $Shape.Top = ($Top / 12700.0) + $halfEmuPt
    $Shape.Width = ($Width / 12700.0) + $halfEmuPt
    $Shape.Height = ($Height / 12700.0) + $halfEmuPt
}

$p = $ppt.ActivePresentation

# --- Slide 17: "filter(function, list)" textbox - off.x/ext.cx nudge ---
$s17 = $p.Slides.Item(17)
$sh17 = $s17.Shapes.Item(5)
Set-ShapeEmu $sh17 3751341 2689393 5502117 774572

# --- Slide 18: ">>>zip([1, 2, 3], [4, 5, 6])..." textbox - off.y nudge ---
$s18 = $p.Slides.Item(18)
$sh18 = $s18.Shapes.Item(7)
Set-ShapeEmu $sh18 3763326 4532865 5478148 1420304

# --- Slide 20: fix wrong result text "210" -> "[1, 2, 3, 4, 5]" ---
$s20 = $p.Slides.Item(20)
$sh20 = $s20.Shapes.Item(5)
$resultParagraph = $sh20.TextFrame.TextRange.Paragraphs(2, 1)
$resultParagraph.Text = "[1, 2, 3, 4, 5]"

# --- Slide 3: "Примеры:" textbox - ext.cy nudge ---
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(5)
Set-ShapeEmu $sh3 909307 1550378 2929948 972254

# --- Slide 8: "for val in x:..." textbox - off.y/ext.cy nudge ---
$s8 = $p.Slides.Item(8)
$sh8 = $s8.Shapes.Item(4)
Set-ShapeEmu $sh8 5052921 1973398 3939028 1336725
